$wb = $excel.ActiveWorkbook

$wsAlerts = $wb.Worksheets.Item('ALERTS')
$wsAlerts.Range('A8:F8').NumberFormat = '@'
$wsAlerts.Range('A8').Value = '2026-01-28'
$wsAlerts.Range('B8').Value = '14:56:46'
$wsAlerts.Range('C8').Value = '14:00'
$wsAlerts.Range('D8').Value = 'Bathroom'
$wsAlerts.Range('E8').Value = 'CRITICAL'
$wsAlerts.Range('F8').Value = 'CRITICAL ALERT: Bathroom occupied, no motion > 60s.'

$wsPir = $wb.Worksheets.Item('PIR')
$wsPir.Range('A110:F122').NumberFormat = '@'
$wsPir.Range('A110').Value = '2026-01-28'
$wsPir.Range('B110').Value = '14:56:39'
$wsPir.Range('C110').Value = '14:00'
$wsPir.Range('D110').Value = 'Bathroom'
$wsPir.Range('E110').Value = 'No Motion'
$wsPir.Range('F110').Value = 'Inactive'
$wsPir.Range('A111').Value = '2026-01-28'
$wsPir.Range('B111').Value = '14:56:43'
$wsPir.Range('C111').Value = '14:00'
$wsPir.Range('D111').Value = 'Bathroom'
$wsPir.Range('E111').Value = 'No Motion'
$wsPir.Range('F111').Value = 'Inactive'
$wsPir.Range('A112').Value = '2026-01-28'
$wsPir.Range('B112').Value = '14:56:48'
$wsPir.Range('C112').Value = '14:00'
$wsPir.Range('D112').Value = 'Bathroom'
$wsPir.Range('E112').Value = 'No Motion'
$wsPir.Range('F112').Value = 'Inactive'
$wsPir.Range('A113').Value = '2026-01-28'
$wsPir.Range('B113').Value = '14:56:54'
$wsPir.Range('C113').Value = '14:00'
$wsPir.Range('D113').Value = 'Bathroom'
$wsPir.Range('E113').Value = 'No Motion'
$wsPir.Range('F113').Value = 'Inactive'
$wsPir.Range('A114').Value = '2026-01-28'
$wsPir.Range('B114').Value = '14:56:58'
$wsPir.Range('C114').Value = '14:00'
$wsPir.Range('D114').Value = 'Bathroom'
$wsPir.Range('E114').Value = 'No Motion'
$wsPir.Range('F114').Value = 'Inactive'
$wsPir.Range('A115').Value = '2026-01-28'
$wsPir.Range('B115').Value = '14:57:03'
$wsPir.Range('C115').Value = '14:00'
$wsPir.Range('D115').Value = 'Bathroom'
$wsPir.Range('E115').Value = 'No Motion'
$wsPir.Range('F115').Value = 'Inactive'
$wsPir.Range('A116').Value = '2026-01-28'
$wsPir.Range('B116').Value = '14:57:08'
$wsPir.Range('C116').Value = '14:00'
$wsPir.Range('D116').Value = 'Bathroom'
$wsPir.Range('E116').Value = 'No Motion'
$wsPir.Range('F116').Value = 'Inactive'
$wsPir.Range('A117').Value = '2026-01-28'
$wsPir.Range('B117').Value = '14:57:14'
$wsPir.Range('C117').Value = '14:00'
$wsPir.Range('D117').Value = 'Bathroom'
$wsPir.Range('E117').Value = 'No Motion'
$wsPir.Range('F117').Value = 'Inactive'
$wsPir.Range('A118').Value = '2026-01-28'
$wsPir.Range('B118').Value = '14:57:18'
$wsPir.Range('C118').Value = '14:00'
$wsPir.Range('D118').Value = 'Bathroom'
$wsPir.Range('E118').Value = 'No Motion'
$wsPir.Range('F118').Value = 'Inactive'
$wsPir.Range('A119').Value = '2026-01-28'
$wsPir.Range('B119').Value = '14:57:24'
$wsPir.Range('C119').Value = '14:00'
$wsPir.Range('D119').Value = 'Bathroom'
$wsPir.Range('E119').Value = 'No Motion'
$wsPir.Range('F119').Value = 'Inactive'
$wsPir.Range('A120').Value = '2026-01-28'
$wsPir.Range('B120').Value = '14:57:29'
$wsPir.Range('C120').Value = '14:00'
$wsPir.Range('D120').Value = 'Bathroom'
$wsPir.Range('E120').Value = 'No Motion'
$wsPir.Range('F120').Value = 'Inactive'
$wsPir.Range('A121').Value = '2026-01-28'
$wsPir.Range('B121').Value = '14:57:34'
$wsPir.Range('C121').Value = '14:00'
$wsPir.Range('D121').Value = 'Bathroom'
$wsPir.Range('E121').Value = 'No Motion'
$wsPir.Range('F121').Value = 'Inactive'
$wsPir.Range('A122').Value = '2026-01-28'
$wsPir.Range('B122').Value = '14:57:39'
$wsPir.Range('C122').Value = '14:00'
$wsPir.Range('D122').Value = 'Bathroom'
$wsPir.Range('E122').Value = 'No Motion'
$wsPir.Range('F122').Value = 'Inactive'

$wsHumidity = $wb.Worksheets.Item('Humidity')
$wsHumidity.Range('A104:F116').NumberFormat = '@'
$wsHumidity.Range('A104').Value = '2026-01-28'
$wsHumidity.Range('B104').Value = '14:56:39'
$wsHumidity.Range('C104').Value = '14:00'
$wsHumidity.Range('D104').Value = 'Bathroom'
$wsHumidity.Range('E104').Value = '88.1%'
$wsHumidity.Range('F104').Value = 'Active'
$wsHumidity.Range('A105').Value = '2026-01-28'
$wsHumidity.Range('B105').Value = '14:56:45'
$wsHumidity.Range('C105').Value = '14:00'
$wsHumidity.Range('D105').Value = 'Bathroom'
$wsHumidity.Range('E105').Value = '88.1%'
$wsHumidity.Range('F105').Value = 'Active'
$wsHumidity.Range('A106').Value = '2026-01-28'
$wsHumidity.Range('B106').Value = '14:56:53'
$wsHumidity.Range('C106').Value = '14:00'
$wsHumidity.Range('D106').Value = 'Bathroom'
$wsHumidity.Range('E106').Value = '87.1%'
$wsHumidity.Range('F106').Value = 'Active'
$wsHumidity.Range('A107').Value = '2026-01-28'
$wsHumidity.Range('B107').Value = '14:57:01'
$wsHumidity.Range('C107').Value = '14:00'
$wsHumidity.Range('D107').Value = 'Bathroom'
$wsHumidity.Range('E107').Value = '87.1%'
$wsHumidity.Range('F107').Value = 'Active'
$wsHumidity.Range('A108').Value = '2026-01-28'
$wsHumidity.Range('B108').Value = '14:57:05'
$wsHumidity.Range('C108').Value = '14:00'
$wsHumidity.Range('D108').Value = 'Bathroom'
$wsHumidity.Range('E108').Value = '87.9%'
$wsHumidity.Range('F108').Value = 'Active'
$wsHumidity.Range('A109').Value = '2026-01-28'
$wsHumidity.Range('B109').Value = '14:57:09'
$wsHumidity.Range('C109').Value = '14:00'
$wsHumidity.Range('D109').Value = 'Bathroom'
$wsHumidity.Range('E109').Value = '87.9%'
$wsHumidity.Range('F109').Value = 'Active'
$wsHumidity.Range('A110').Value = '2026-01-28'
$wsHumidity.Range('B110').Value = '14:57:13'
$wsHumidity.Range('C110').Value = '14:00'
$wsHumidity.Range('D110').Value = 'Bathroom'
$wsHumidity.Range('E110').Value = '87.9%'
$wsHumidity.Range('F110').Value = 'Active'
$wsHumidity.Range('A111').Value = '2026-01-28'
$wsHumidity.Range('B111').Value = '14:57:17'
$wsHumidity.Range('C111').Value = '14:00'
$wsHumidity.Range('D111').Value = 'Bathroom'
$wsHumidity.Range('E111').Value = '87.9%'
$wsHumidity.Range('F111').Value = 'Active'
$wsHumidity.Range('A112').Value = '2026-01-28'
$wsHumidity.Range('B112').Value = '14:57:21'
$wsHumidity.Range('C112').Value = '14:00'
$wsHumidity.Range('D112').Value = 'Bathroom'
$wsHumidity.Range('E112').Value = '87.0%'
$wsHumidity.Range('F112').Value = 'Active'
$wsHumidity.Range('A113').Value = '2026-01-28'
$wsHumidity.Range('B113').Value = '14:57:25'
$wsHumidity.Range('C113').Value = '14:00'
$wsHumidity.Range('D113').Value = 'Bathroom'
$wsHumidity.Range('E113').Value = '87.9%'
$wsHumidity.Range('F113').Value = 'Active'
$wsHumidity.Range('A114').Value = '2026-01-28'
$wsHumidity.Range('B114').Value = '14:57:30'
$wsHumidity.Range('C114').Value = '14:00'
$wsHumidity.Range('D114').Value = 'Bathroom'
$wsHumidity.Range('E114').Value = '87.9%'
$wsHumidity.Range('F114').Value = 'Active'
$wsHumidity.Range('A115').Value = '2026-01-28'
$wsHumidity.Range('B115').Value = '14:57:33'
$wsHumidity.Range('C115').Value = '14:00'
$wsHumidity.Range('D115').Value = 'Bathroom'
$wsHumidity.Range('E115').Value = '87.0%'
$wsHumidity.Range('F115').Value = 'Active'
$wsHumidity.Range('A116').Value = '2026-01-28'
$wsHumidity.Range('B116').Value = '14:57:37'
$wsHumidity.Range('C116').Value = '14:00'
$wsHumidity.Range('D116').Value = 'Bathroom'
$wsHumidity.Range('E116').Value = '87.9%'
$wsHumidity.Range('F116').Value = 'Active'

$wsTemperature = $wb.Worksheets.Item('Temperature')
$wsTemperature.Range('A104:F116').NumberFormat = '@'
$wsTemperature.Range('A104').Value = '2026-01-28'
$wsTemperature.Range('B104').Value = '14:56:39'
$wsTemperature.Range('C104').Value = '14:00'
$wsTemperature.Range('D104').Value = 'Bathroom'
$wsTemperature.Range('E104').Value = '23.0C'
$wsTemperature.Range('F104').Value = 'Active'
$wsTemperature.Range('A105').Value = '2026-01-28'
$wsTemperature.Range('B105').Value = '14:56:45'
$wsTemperature.Range('C105').Value = '14:00'
$wsTemperature.Range('D105').Value = 'Bathroom'
$wsTemperature.Range('E105').Value = '23.0C'
$wsTemperature.Range('F105').Value = 'Active'
$wsTemperature.Range('A106').Value = '2026-01-28'
$wsTemperature.Range('B106').Value = '14:56:53'
$wsTemperature.Range('C106').Value = '14:00'
$wsTemperature.Range('D106').Value = 'Bathroom'
$wsTemperature.Range('E106').Value = '23.0C'
$wsTemperature.Range('F106').Value = 'Active'
$wsTemperature.Range('A107').Value = '2026-01-28'
$wsTemperature.Range('B107').Value = '14:57:01'
$wsTemperature.Range('C107').Value = '14:00'
$wsTemperature.Range('D107').Value = 'Bathroom'
$wsTemperature.Range('E107').Value = '23.0C'
$wsTemperature.Range('F107').Value = 'Active'
$wsTemperature.Range('A108').Value = '2026-01-28'
$wsTemperature.Range('B108').Value = '14:57:05'
$wsTemperature.Range('C108').Value = '14:00'
$wsTemperature.Range('D108').Value = 'Bathroom'
$wsTemperature.Range('E108').Value = '22.9C'
$wsTemperature.Range('F108').Value = 'Active'
$wsTemperature.Range('A109').Value = '2026-01-28'
$wsTemperature.Range('B109').Value = '14:57:09'
$wsTemperature.Range('C109').Value = '14:00'
$wsTemperature.Range('D109').Value = 'Bathroom'
$wsTemperature.Range('E109').Value = '23.0C'
$wsTemperature.Range('F109').Value = 'Active'
$wsTemperature.Range('A110').Value = '2026-01-28'
$wsTemperature.Range('B110').Value = '14:57:13'
$wsTemperature.Range('C110').Value = '14:00'
$wsTemperature.Range('D110').Value = 'Bathroom'
$wsTemperature.Range('E110').Value = '22.9C'
$wsTemperature.Range('F110').Value = 'Active'
$wsTemperature.Range('A111').Value = '2026-01-28'
$wsTemperature.Range('B111').Value = '14:57:17'
$wsTemperature.Range('C111').Value = '14:00'
$wsTemperature.Range('D111').Value = 'Bathroom'
$wsTemperature.Range('E111').Value = '22.9C'
$wsTemperature.Range('F111').Value = 'Active'
$wsTemperature.Range('A112').Value = '2026-01-28'
$wsTemperature.Range('B112').Value = '14:57:21'
$wsTemperature.Range('C112').Value = '14:00'
$wsTemperature.Range('D112').Value = 'Bathroom'
$wsTemperature.Range('E112').Value = '22.9C'
$wsTemperature.Range('F112').Value = 'Active'
$wsTemperature.Range('A113').Value = '2026-01-28'
$wsTemperature.Range('B113').Value = '14:57:26'
$wsTemperature.Range('C113').Value = '14:00'
$wsTemperature.Range('D113').Value = 'Bathroom'
$wsTemperature.Range('E113').Value = '22.9C'
$wsTemperature.Range('F113').Value = 'Active'
$wsTemperature.Range('A114').Value = '2026-01-28'
$wsTemperature.Range('B114').Value = '14:57:30'
$wsTemperature.Range('C114').Value = '14:00'
$wsTemperature.Range('D114').Value = 'Bathroom'
$wsTemperature.Range('E114').Value = '22.9C'
$wsTemperature.Range('F114').Value = 'Active'
$wsTemperature.Range('A115').Value = '2026-01-28'
$wsTemperature.Range('B115').Value = '14:57:34'
$wsTemperature.Range('C115').Value = '14:00'
$wsTemperature.Range('D115').Value = 'Bathroom'
$wsTemperature.Range('E115').Value = '22.9C'
$wsTemperature.Range('F115').Value = 'Active'
$wsTemperature.Range('A116').Value = '2026-01-28'
$wsTemperature.Range('B116').Value = '14:57:38'
$wsTemperature.Range('C116').Value = '14:00'
$wsTemperature.Range('D116').Value = 'Bathroom'
$wsTemperature.Range('E116').Value = '22.9C'
$wsTemperature.Range('F116').Value = 'Active'

$wsMmwave = $wb.Worksheets.Item('mmWave')
$wsMmwave.Range('A5:F5').NumberFormat = '@'
$wsMmwave.Range('A5').Value = '2026-01-28'
$wsMmwave.Range('B5').Value = '14:57:22'
$wsMmwave.Range('C5').Value = '14:00'
$wsMmwave.Range('D5').Value = 'Living Room'
$wsMmwave.Range('E5').Value = 'Presence Detected'
$wsMmwave.Range('F5').Value = 'Active'

Write-Output 'Edit complete'